$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.986.41"
$ws.Range("E2").Value = "  -1.47%  "

$ws.Range("D3").Value = "2.009.36"
$ws.Range("E3").Value = "  -2.91%  "

$ws.Range("E4").Value = "  -0.24%  "

$ws.Range("D5").Value = "'224.97"
$ws.Range("E5").Value = "  -2.73%  "

$ws.Range("D6").Value = "'0.605"
$ws.Range("E6").Value = "  -2.67%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "'54.74"
$ws.Range("E8").Value = "  -5.35%  "

$ws.Range("E9").Value = "  -2.93%  "

$ws.Range("D10").Value = "'0.0781"
$ws.Range("E10").Value = "  +0.90%  "

$ws.Range("E11").Value = "  -4.85%  "

$ws.Range("D12").Value = "2.306.51"
$ws.Range("E12").Value = "  -2.95%  "

$ws.Range("E13").Value = "  -4.54%  "

$ws.Range("D14").Value = "'20.14"
$ws.Range("E14").Value = "  -4.72%  "

$ws.Range("E15").Value = "  -3.33%  "

$ws.Range("E16").Value = "  -4.07%  "

$ws.Range("D17").Value = "2.011.10"
$ws.Range("E17").Value = "  -3.07%  "

$ws.Range("D18").Value = "36.962.99"
$ws.Range("E18").Value = "  -1.50%  "

$ws.Range("D19").Value = "'6.17"
$ws.Range("E19").Value = "  -0.16%  "

$ws.Range("D20").Value = "'68.59"
$ws.Range("E20").Value = "  -1.85%  "

$ws.Range("D21").Value = "0.0₃0811"
$ws.Range("E21").Value = "  -1.64%  "

$ws.Range("D22").Value = "'222.61"
$ws.Range("E22").Value = "  -1.90%  "

$ws.Range("E23").Value = "  -0.01%  "

$ws.Range("E24").Value = "  +0.96%  "

$ws.Range("D25").Value = "'2.17"
$ws.Range("E25").Value = "  -6.76%  "

$ws.Range("D26").Value = "'165.91"
$ws.Range("E26").Value = "  -2.07%  "

$ws.Range("E27").Value = "  -7.68%  "

$ws.Range("E28").Value = "  -0.04%  "

$ws.Range("E29").Value = "  -3.37%  "

$ws.Range("E30").Value = "  -5.08%  "

$ws.Range("E31").Value = "  -4.05%  "

$ws.Range("D32").Value = "'4.50"
$ws.Range("E32").Value = "  -1.37%  "

$ws.Range("E33").Value = "  -2.39%  "

$ws.Range("D34").Value = "'4.39"
$ws.Range("E34").Value = "  -5.10%  "

$ws.Range("D35").Value = "'2.33"
$ws.Range("E35").Value = "  -7.89%  "

$ws.Range("D36").Value = "'1.84"
$ws.Range("E36").Value = "  +1.10%  "

$ws.Range("E37").Value = "  -0.19%  "

$ws.Range("E38").Value = "  -4.43%  "

$ws.Range("E39").Value = "  -1.16%  "

$ws.Range("D40").Value = "1.478.72"
$ws.Range("E40").Value = "  -0.48%  "

$ws.Range("E41").Value = "  -5.23%  "

$ws.Range("D42").Value = "'94.67"
$ws.Range("E42").Value = "  -3.57%  "

$ws.Range("D43").Value = "'0.0915"
$ws.Range("E43").Value = "  -4.49%  "

$ws.Range("D44").Value = "'16.27"
$ws.Range("E44").Value = "  -2.24%  "

$ws.Range("D45").Value = "'2.76"
$ws.Range("E45").Value = "  -5.14%  "

$ws.Range("E46").Value = "  -5.84%  "

$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'7.18"
$ws.Range("E47").Value = "  -0.82%  "

$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "'1.01"
$ws.Range("E48").Value = "  -2.73%  "

$ws.Range("E49").Value = "  -1.20%  "

$ws.Range("D50").Value = "2.193.15"
$ws.Range("E50").Value = "  -3.00%  "

$ws.Range("E51").Value = "  -3.30%  "

